$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.2609752687857849, 0.9909744767902922, 2.68705854512723, 1.639224983071949, 1.637933849261315, 42)
    3 = @(0.009041043032269717, 1.032344038245651, 2.857553663655543, 1.690430023294529, 1.711405481256953, 41)
    4 = @(0.2438441126156719, 1.021400937554764, 2.863247739468819, 1.692113394388455, 1.695782940939382, 40)
    5 = @(0.05021881217214593, 1.021229807756068, 2.921706329422886, 1.709299953028399, 1.730897178471549, 39)
    6 = @(0.2407774045323826, 1.05914724455951, 3.019659849071489, 1.737716849510152, 1.744056094439822, 38)
    7 = @(0.08256934594935905, 1.017230331921043, 2.994009761198107, 1.730320710503723, 1.752189950219291, 37)
    8 = @(0.2581920007548112, 1.070846802007395, 3.103390552560661, 1.761644275261229, 1.767340115694834, 36)
    9 = @(0.05684027587393566, 0.9604249871766081, 2.945262510222375, 1.71617671299385, 1.740276425508428, 35)
    10 = @(0.2309292503718555, 1.011888791862958, 3.109978691686091, 1.763513167426342, 1.774619999292197, 34)
    11 = @(0.1436658260302781, 0.9494617496012048, 3.096262616726478, 1.759620020551732, 1.780936820329735, 33)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = $vals[0]
    $ws.Cells.Item($rowNum, 3).Value = $vals[1]
    $ws.Cells.Item($rowNum, 4).Value = $vals[2]
    $ws.Cells.Item($rowNum, 5).Value = $vals[3]
    $ws.Cells.Item($rowNum, 6).Value = $vals[4]
    $ws.Cells.Item($rowNum, 7).Value = $vals[5]
}
